# Update "想去人数" (F column, number of people interested) figures that were
# refreshed by the site generator, plus one ticket-status flip (G12: a
# numeric lowest price -> "不可售" / "not for sale") on the sheets that carry
# that same event row.
#
# Sheet map (by tab order):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life) - no data rows, untouched
#   4 = 全部类型  (All types, a merge of sheet 1 + sheet 2)

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 -------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 6793
$ws1.Range("F3").Value  = 87
$ws1.Range("F5").Value  = 438
$ws1.Range("F6").Value  = 144
$ws1.Range("F7").Value  = 6415
$ws1.Range("F8").Value  = 55
$ws1.Range("F9").Value  = 194
$ws1.Range("F12").Value = 106
$ws1.Range("G12").Value = "不可售"
$ws1.Range("F13").Value = 392
$ws1.Range("F14").Value = 129
$ws1.Range("F16").Value = 373
$ws1.Range("F17").Value = 45
$ws1.Range("F18").Value = 8
$ws1.Range("F19").Value = 4761
$ws1.Range("F20").Value = 76
$ws1.Range("F21").Value = 55
$ws1.Range("F22").Value = 218
$ws1.Range("F23").Value = 203
$ws1.Range("F24").Value = 113

# ---- Sheet 2: 演出 --------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F2").Value = 44

# ---- Sheet 4: 全部类型 ----------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 6793
$ws4.Range("F4").Value  = 18
$ws4.Range("F5").Value  = 438
$ws4.Range("F7").Value  = 6415
$ws4.Range("F8").Value  = 55
$ws4.Range("F9").Value  = 194
$ws4.Range("F10").Value = 1275
$ws4.Range("F11").Value = 11
$ws4.Range("F12").Value = 106
$ws4.Range("G12").Value = "不可售"
$ws4.Range("F14").Value = 129
$ws4.Range("F16").Value = 373
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 4761
$ws4.Range("F20").Value = 44
$ws4.Range("F21").Value = 76
$ws4.Range("F23").Value = 218
$ws4.Range("F24").Value = 203
$ws4.Range("F25").Value = 113
